# Donate button group (slide 4, "Group 6") is flattened into its two
# member shapes (the rounded-rectangle "Donate" caption and the arrow
# picture) which become direct children of the slide, matching the
# capitalization fix commit that removed the wrapping group.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# Locate the group shape by name (robust to ordering) and ungroup it.
$grp = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Group 6") {
        $grp = $sh
        break
    }
}

if ($grp -ne $null) {
    $grp.Ungroup() | Out-Null
}
